$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename readout labels from "Weight ..." to "mass of ..."
$ws.Range("E2").Value = " mass of 2MIM_2"
$ws.Range("E4").Value = " mass of MeOH_2"

# Rename reagent IDs from "...002" lot numbers to "..._2" series names
$ws.Range("F3").Value = " 2MIM_2"
$ws.Range("F5").Value = " MeOH_2"

# Update the active selection to F6
$ws.Range("F6").Select()
